# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 16:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 561159
$ws.Range("C4").Value = 859
$ws.Range("D4").Value = 33122
$ws.Range("E4").Value = 505904
$ws.Range("F4").Value = 11770
$ws.Range("G4").Value = 28
$ws.Range("H4").Value = 22133

# Row 15 - Suiza
$ws.Range("B15").Value = 25623
$ws.Range("C15").Value = 208
$ws.Range("E15").Value = 11794
$ws.Range("G15").Value = 23
$ws.Range("H15").Value = 1129

# Row 20 - Austria
$ws.Range("B20").Value = 14013
$ws.Range("C20").Value = 68
$ws.Range("E20").Value = 6302

# Row 54 - Argentina
$ws.Range("E54").Value = 1644
$ws.Range("G54").Value = 6
$ws.Range("H54").Value = 96

# Row 70 - Azerbaiyan
$ws.Range("B70").Value = 1148
$ws.Range("C70").Value = 50
$ws.Range("D70").Value = 289
$ws.Range("E70").Value = 847
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 12

# Row 88 - Principado de Andorra
$ws.Range("B88").Value = 646
$ws.Range("C88").Value = 8
$ws.Range("E88").Value = 489

# Row 118 - Kenia
$ws.Range("B118").Value = 208
$ws.Range("C118").Value = 11
$ws.Range("D118").Value = 40
$ws.Range("E118").Value = 159
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 9
